$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A labels (updated) and Column B values (updated) for rows 1-43
$data = @(
    @(1,  "Total time taken for the ride", 0.03384300925925926),
    @(2,  "Actual Ampere-hours (Ah)", 27.74537166666667),
    @(3,  "Actual Watt-hours (Wh)", 1385.207483403889),
    @(4,  "Starting SoC (Ah)", 38.745),
    @(5,  "Ending SoC (Ah)", 10.26),
    @(6,  "Starting SoC (%)", 97),
    @(7,  "Ending SoC (%)", 25),
    @(8,  "Total distance covered (km)", 27.14508620657145),
    @(9,  "Total energy consumption(WH/KM)", 51.02976917673407),
    @(10, "Total SOC consumed(%)", 72),
    @(11, "Mode", "Custom mode`n92.78%`nEco mode`n6.14%`nSports mode`n0.07%"),
    @(12, "Peak Power(kW)", 5724.067349),
    @(13, "Average Power(kW)", -1717.199359591598),
    @(14, "Total Energy Regenerated(kWh)", 0.02443286916666667),
    @(15, "Regenerative Effectiveness(%)", 0.001763810729426659),
    @(16, "Highest Cell Voltage(V)", 3.326),
    @(17, "Lowest Cell Voltage(V)", 2.963),
    @(18, "Difference in Cell Voltage(V)", 0.363),
    @(19, "Minimum Temperature(C)", 25),
    @(20, "Maximum Temperature(C)", 40),
    @(21, "Difference in Temperature(C)", 15),
    @(22, "Maximum Fet Temperature-BMS(C)", 70),
    @(23, "Maximum Afe Temperature-BMS(C)", 63),
    @(24, "Maximum PCB Temperature-BMS(C)", 64),
    @(25, "Maximum MCU Temperature(C)", 60),
    @(26, "Maximum Motor Temperature(C)", 95),
    @(27, "Abnormal Motor Temperature Detected(C)", 0),
    @(28, "highest cell temp(C)", 40),
    @(29, "lowest cell temp(C)", 25),
    @(30, "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)", 15),
    @(31, "Battery Voltage(V)", 55),
    @(32, "Total energy charged(kWh)", 1.525995441666667),
    @(33, "Electricity consumption units(kW)", [double]"1.449684072110757e-07"),
    @(34, "Idling time percentage", 23.11878333573591),
    @(35, "Time spent in 0-10 km/h", 8.433040219114892),
    @(36, "Time spent in 10-20 km/h", 6.804093988755946),
    @(37, "Time spent in 20-30 km/h", 7.737494594204988),
    @(38, "Time spent in 30-40 km/h", 8.609629522848493),
    @(39, "Time spent in 40-50 km/h", 11.99365720051896),
    @(40, "Time spent in 50-60 km/h", 8.353755225601846),
    @(41, "Time spent in 60-70 km/h", 13.13968574311662),
    @(42, "Time spent in 70-80 km/h", 11.66210177310076),
    @(43, "Time spent in 80-90 km/h", 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $label = $row[1]
    $value = $row[2]
    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $value
}
